$d = $word.ActiveDocument

# Step 1: Move the "_GoBack" bookmark from the final (empty) paragraph to
# immediately after the text run in the first paragraph (before its
# paragraph mark). Bookmark names are unique, so re-adding "_GoBack" at the
# new location removes it from its old location automatically.
$findRng = $d.Content
$found = $findRng.Find.Execute("minimizes error {the difference between the true measurements and their associated matches generated from the map using the guessed pose.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target paragraph 1 text to anchor the _GoBack bookmark."
}
$d.Bookmarks.Add("_GoBack", $findRng)

# Step 2: Replace the contents of the now-empty last paragraph (which used
# to hold the bookmark) with the new "LSD SLAM" discussion, and add a new
# paragraph after it with the "ORB SLAM" discussion.
$lastParaIndex = $d.Paragraphs.Count
$targetRange = $d.Paragraphs.Item($lastParaIndex).Range

$insertXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">The tracking thread in </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>LSD SLAM</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> is responsible for estimating the pose of the current frame with respect to the currently active keyframe in the map, using the previous frame pose as a prior. The required pose is represented by an </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>SE(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>3) transformation, and is found by an iteratively re-weighted Gauss-Newton optimization that minimizes the variance normalized photometric residual error, as described in [78]. A keyframe is considered active if it is the most recent keyframe accommodated in the map. To minimize outlier effects, measurements with large residuals are down-weighted from one iteration to the next.</w:t>
      </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">Pose estimation in </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>ORB SLAM</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> is established through a constant velocity motion model prior, followed by a pose reﬁnement using optimization. As the motion model is expected to be easily violated through abrupt motions, ORB SLAM detects such failures by tracking the number of matched features; if it falls below a certain threshold, map points are projected onto the current frame, and a wide-range feature search takes place around the projected locations. </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>In an effort to</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> make ORB SLAM operate in large environments, a subset of the global map, known as the local map, is deﬁned by all landmarks corresponding to the set of all keyframes that share edges with the current frame, as well as all neighbors of this set of keyframes from the pose graph. The selected landmarks are ﬁltered out to keep </w:t>
      </w:r>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>only the features that are most likely to be matched in the current frame. Furthermore, if the distance from the camera’s center to the landmark is beyond the range of the valid features, the landmark is also discarded. The remaining set of landmarks is then searched for and matched in the current frame, before a ﬁnal camera pose reﬁnement step.</w:t>
      </w:r></w:p>

'@

$targetRange.InsertXML($insertXml)
